$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.538445
$ws.Range("H2").Value = 1.615335
$ws.Range("I2").Value = 0.03371608002174246
$ws.Range("J2").Value = 0.03371608002174246
$ws.Range("M2").Value = 0.2799683333333333
$ws.Range("N2").Value = 0.839905
$ws.Range("O2").Value = 0.0294305463214559
$ws.Range("P2").Value = 0.0294305463214559
$ws.Range("Q2").Value = 0.1507475492416666
$ws.Range("R2").Value = 1.356727943175
$ws.Range("S2").Value = 0.0009922826548578055
$ws.Range("T2").Value = 0.0009922826548578055
$ws.Range("G3").Value = 0.538445
$ws.Range("H3").Value = 1.615335
$ws.Range("I3").Value = 0.03371608002174246
$ws.Range("J3").Value = 0.03371608002174246
$ws.Range("O3").Value = 0.2486942046732164
$ws.Range("P3").Value = 0.2486942046732163
$ws.Range("Q3").Value = 1.273847976031667
$ws.Range("R3").Value = 11.464631784285
$ws.Range("S3").Value = 0.008384993705705761
$ws.Range("T3").Value = 0.008384993705705761
$ws.Range("G4").Value = 0.538445
$ws.Range("H4").Value = 1.615335
$ws.Range("I4").Value = 0.03371608002174246
$ws.Range("J4").Value = 0.03371608002174246
$ws.Range("M4").Value = 6.86709
$ws.Range("N4").Value = 20.60127
$ws.Range("O4").Value = 0.7218752490053277
$ws.Range("P4").Value = 0.7218752490053277
$ws.Range("Q4").Value = 3.69755027505
$ws.Range("R4").Value = 33.27795247545
$ws.Range("S4").Value = 0.0243388036611789
$ws.Range("T4").Value = 0.0243388036611789
$ws.Range("I5").Value = 0.7539416098905094
$ws.Range("J5").Value = 0.7539416098905093
$ws.Range("M5").Value = 0.2799683333333333
$ws.Range("N5").Value = 0.839905
$ws.Range("O5").Value = 0.0294305463214559
$ws.Range("P5").Value = 0.0294305463214559
$ws.Range("Q5").Value = 3.370939026393889
$ws.Range("R5").Value = 30.33845123754501
$ws.Range("S5").Value = 0.02218891347355567
$ws.Range("T5").Value = 0.02218891347355567
$ws.Range("I6").Value = 0.7539416098905094
$ws.Range("J6").Value = 0.7539416098905093
$ws.Range("O6").Value = 0.2486942046732164
$ws.Range("P6").Value = 0.2486942046732163
$ws.Range("R6").Value = 256.3661890312191
$ws.Range("S6").Value = 0.1875009090417646
$ws.Range("T6").Value = 0.1875009090417646
$ws.Range("I7").Value = 0.7539416098905094
$ws.Range("J7").Value = 0.7539416098905093
$ws.Range("M7").Value = 6.86709
$ws.Range("N7").Value = 20.60127
$ws.Range("O7").Value = 0.7218752490053277
$ws.Range("P7").Value = 0.7218752490053277
$ws.Range("Q7").Value = 82.68271415967001
$ws.Range("R7").Value = 744.1444274370301
$ws.Range("S7").Value = 0.5442517873751891
$ws.Range("T7").Value = 0.5442517873751891
$ws.Range("G8").Value = 3.391101666666666
$ws.Range("H8").Value = 10.173305
$ws.Range("I8").Value = 0.2123423100877482
$ws.Range("J8").Value = 0.2123423100877481
$ws.Range("M8").Value = 0.2799683333333333
$ws.Range("N8").Value = 0.839905
$ws.Range("O8").Value = 0.0294305463214559
$ws.Range("P8").Value = 0.0294305463214559
$ws.Range("Q8").Value = 0.9494010817805554
$ws.Range("R8").Value = 8.544609736024999
$ws.Range("S8").Value = 0.006249350193042425
$ws.Range("T8").Value = 0.006249350193042425
$ws.Range("G9").Value = 3.391101666666666
$ws.Range("H9").Value = 10.173305
$ws.Range("I9").Value = 0.2123423100877482
$ws.Range("J9").Value = 0.2123423100877481
$ws.Range("O9").Value = 0.2486942046732164
$ws.Range("P9").Value = 0.2486942046732163
$ws.Range("Q9").Value = 8.022635542350555
$ws.Range("R9").Value = 72.203719881155
$ws.Range("S9").Value = 0.05280830192574602
$ws.Range("T9").Value = 0.05280830192574601
$ws.Range("G10").Value = 3.391101666666666
$ws.Range("H10").Value = 10.173305
$ws.Range("I10").Value = 0.2123423100877482
$ws.Range("J10").Value = 0.2123423100877481
$ws.Range("M10").Value = 6.86709
$ws.Range("N10").Value = 20.60127
$ws.Range("O10").Value = 0.7218752490053277
$ws.Range("P10").Value = 0.7218752490053277
$ws.Range("Q10").Value = 23.28700034415
$ws.Range("R10").Value = 209.58300309735
$ws.Range("S10").Value = 0.1532846579689597
$ws.Range("T10").Value = 0.1532846579689597
